$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 35742
$ws.Range("J3").Value = 35742
$ws.Range("L3").Value = 35742
$ws.Range("N3").Value = -35970
$ws.Range("H9").Value = 2053.5454
$ws.Range("I9").Value = 441.2857
$ws.Range("K9").Value = 441.2857
$ws.Range("M9").Value = -272.2857
$ws.Range("H11").Value = 4954.8335
$ws.Range("I11").Value = 4954.8335
$ws.Range("K11").Value = 4954.8335
$ws.Range("M11").Value = -4814.8335
$ws.Range("H29").Value = 49
$ws.Range("I29").Value = 49
$ws.Range("K29").Value = 147
$ws.Range("M29").Value = 134
$ws.Range("H62").Value = 1996.3334
$ws.Range("I62").Value = 1996.3334
$ws.Range("K62").Value = 1996.3334
$ws.Range("M62").Value = -1372.3334
$ws.Range("H65").Value = 1996.3334
$ws.Range("I65").Value = 1996.3334
$ws.Range("K65").Value = 9981.666999999999
$ws.Range("M65").Value = -6861.666999999999
$ws.Range("H93").Value = 58533.332
$ws.Range("J93").Value = 58533.332
$ws.Range("L93").Value = 58533.332
$ws.Range("N93").Value = -63525.332
$ws.Range("H102").Value = 35742
$ws.Range("J102").Value = 35742
$ws.Range("L102").Value = 35742
$ws.Range("N102").Value = -42232
$ws.Range("H111").Value = 2182.9375
$ws.Range("I111").Value = 2275.8
$ws.Range("J111").Value = 2028.1666
$ws.Range("K111").Value = 6827.400000000001
$ws.Range("L111").Value = 6084.4998
$ws.Range("M111").Value = -3760.400000000001
$ws.Range("N111").Value = -12218.4998
$ws.Range("H137").Value = 2828.1667
$ws.Range("J137").Value = 3247.5
$ws.Range("L137").Value = 9742.5
$ws.Range("N137").Value = -14842.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15439.833
$ws.Range("I32").Value = 15439.833
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 15439.833
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -15152.833
$ws.Range("N32").Value = $null
$ws.Range("H61").Value = 6348.125
$ws.Range("I61").Value = 5357.4
$ws.Range("J61").Value = 7999.3335
$ws.Range("K61").Value = 5357.4
$ws.Range("L61").Value = 7999.3335
$ws.Range("M61").Value = -5145.4
$ws.Range("N61").Value = -8423.333500000001
$ws.Range("H74").Value = 2399
$ws.Range("I74").Value = 2399
$ws.Range("K74").Value = 2399
$ws.Range("M74").Value = -1525
$ws.Range("H77").Value = 2399
$ws.Range("I77").Value = 2399
$ws.Range("K77").Value = 11995
$ws.Range("M77").Value = -7627
$ws.Range("H98").Value = 37037.5
$ws.Range("J98").Value = 37037.5
$ws.Range("L98").Value = 37037.5
$ws.Range("N98").Value = -43027.5
$ws.Range("H136").Value = 6348.125
$ws.Range("I136").Value = 5357.4
$ws.Range("J136").Value = 7999.3335
$ws.Range("K136").Value = 16072.2
$ws.Range("L136").Value = 23998.0005
$ws.Range("M136").Value = -13522.2
$ws.Range("N136").Value = -29098.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 8000
$ws.Range("J103").Value = 8000
$ws.Range("L103").Value = 8000
$ws.Range("N103").Value = -10344
$ws.Range("H134").Value = 3219.8
$ws.Range("I134").Value = 1549.5
$ws.Range("J134").Value = 4333.3335
$ws.Range("K134").Value = 4648.5
$ws.Range("L134").Value = 13000.0005
$ws.Range("M134").Value = -2113.5
$ws.Range("N134").Value = -18070.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2419.5454
$ws.Range("I7").Value = 2142.8
$ws.Range("K7").Value = 2142.8
$ws.Range("M7").Value = -2029.8
$ws.Range("H31").Value = 3445.5557
$ws.Range("J31").Value = 4499.5
$ws.Range("L31").Value = 4499.5
$ws.Range("N31").Value = -5089.5
$ws.Range("H34").Value = 3445.5557
$ws.Range("J34").Value = 4499.5
$ws.Range("L34").Value = 4499.5
$ws.Range("N34").Value = -4903.5
$ws.Range("H43").Value = 49999
$ws.Range("J43").Value = 49999
$ws.Range("L43").Value = 49999
$ws.Range("N43").Value = -50367
$ws.Range("H69").Value = 9296.200000000001
$ws.Range("I69").Value = 9296.200000000001
$ws.Range("K69").Value = 9296.200000000001
$ws.Range("M69").Value = -8547.200000000001
$ws.Range("H72").Value = 9296.200000000001
$ws.Range("I72").Value = 9296.200000000001
$ws.Range("K72").Value = 27888.6
$ws.Range("M72").Value = -24144.6
$ws.Range("H101").Value = 49999
$ws.Range("J101").Value = 49999
$ws.Range("L101").Value = 49999
$ws.Range("N101").Value = -56489
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 859.1667
$ws.Range("J12").Value = 909.9091
$ws.Range("L12").Value = 2729.7273
$ws.Range("N12").Value = -3075.7273
$ws.Range("H40").Value = 26
$ws.Range("I40").Value = 50
$ws.Range("K40").Value = 200
$ws.Range("M40").Value = -131
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 964.5
$ws.Range("I2").Value = 20.181818
$ws.Range("J2").Value = 2448.4285
$ws.Range("K2").Value = 20.181818
$ws.Range("L2").Value = 2448.4285
$ws.Range("M2").Value = 92.81818200000001
$ws.Range("N2").Value = -2674.4285
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null
$ws.Range("H100").Value = 100355
$ws.Range("J100").Value = 100355
$ws.Range("L100").Value = 100355
$ws.Range("N100").Value = -102519
$ws.Range("H122").Value = 2890.6956
$ws.Range("I122").Value = 2965.5454
$ws.Range("K122").Value = 8896.636200000001
$ws.Range("M122").Value = -6446.636200000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 769.8333
$ws.Range("I22").Value = 782.25
$ws.Range("J22").Value = 745
$ws.Range("K22").Value = 782.25
$ws.Range("L22").Value = 745
$ws.Range("M22").Value = -487.25
$ws.Range("N22").Value = -1335
$ws.Range("H27").Value = 769.8333
$ws.Range("I27").Value = 782.25
$ws.Range("J27").Value = 745
$ws.Range("K27").Value = 782.25
$ws.Range("L27").Value = 745
$ws.Range("M27").Value = -675.25
$ws.Range("N27").Value = -959
$ws.Range("H103").Value = 15199.4
$ws.Range("J103").Value = 15199.4
$ws.Range("L103").Value = 15199.4
$ws.Range("N103").Value = -17543.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 26067.334
$ws.Range("J101").Value = 26067.334
$ws.Range("L101").Value = 26067.334
$ws.Range("N101").Value = -32557.334
$ws.Range("H103").Value = 35999.8
$ws.Range("J103").Value = 35999.8
$ws.Range("L103").Value = 35999.8
$ws.Range("N103").Value = -38343.8
$ws.Range("H107").Value = 1264.1428
$ws.Range("I107").Value = 1141.8334
$ws.Range("J107").Value = 1998
$ws.Range("K107").Value = 3425.5002
$ws.Range("L107").Value = 5994
$ws.Range("M107").Value = -1505.5002
$ws.Range("N107").Value = -9834
$ws.Range("H136").Value = 48153.727
$ws.Range("I136").Value = 57187.89
$ws.Range("K136").Value = 171563.67
$ws.Range("M136").Value = -169013.67
